$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "LOCATION_ON_04/01/19-18:45"
$ws.Range("B2").Value = "SUBLOCATION_ON_04/01/19-18:45"
$ws.Range("C2").Value = "CUSTOMER_ON_04/01/19-11:48"
$ws.Range("D2").Value = "EQUIPMENT_ON_04/01/19-11:52"
$ws.Range("E2").Value = "ITEM_ON_04/01/19-11:56"
$ws.Range("F2").Value = "SUPPLIER_ON_04/01/19-18:46"
$ws.Range("G2").Value = "SUBCUSTOMER1_ON_04/01/19-11:48"
$ws.Range("H2").Value = "SUBCUSTOMER2_ON_04/01/19-11:48"
$ws.Range("I2").Value = "SUBEQUIPMENT1_ON_04/01/19-11:52"
$ws.Range("J2").Value = "SUBEQUIPMENT2_ON_04/01/19-11:52"
$ws.Range("K2").Value = "SUBITEM1_ON_04/01/19-11:56"
$ws.Range("L2").Value = "SUBITEM2_ON_04/01/19-11:56"
$ws.Range("M2").Value = "SUBSUPPLIER1_ON_04/01/19-18:46"
$ws.Range("N2").Value = "SUBSUPPLIER2_ON_04/01/19-18:46"
$ws.Range("O2").Value = "Customers > Cust11343 > @@@"
$ws.Range("P2").Value = "Suppliers"
$ws.Range("Q2").Value = "SUPPLIER_ON_04/01/19-18:46"
$ws.Range("R2").Value = "SUBSUPPLIER1_ON_04/01/19-18:46"
$ws.Range("S2").Value = "SUBSUPPLIER2_ON_04/01/19-18:46"
